# Generate Report for Handoff
# Swaps the reported rows so that:
#  - d65a6fcc-...md now sits in row 2 (still "Handed back: in sync with en-US")
#  - 4f079f8c-...md moves to row 3 and is now "Ready for handoff" with fresh
#    handoff timestamps, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$fileA = "4f079f8c-42cb-47c3-ad76-af83946074ac"
$fileB = "d65a6fcc-06b9-4407-ba61-b49c531341c6"

$statusHandedBack = "Handed back: in sync with en-US"
$statusReady = "Ready for handoff"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "$fileB.md"
$ws1.Range("B2").Value = $statusHandedBack
$ws1.Range("C2").Value = $statusHandedBack
$ws1.Range("D2").Value = "2016-43-21 04:43:53"

$ws1.Range("A3").Value = "$fileA.md"
$ws1.Range("B3").Value = $statusReady
$ws1.Range("C3").Value = $statusReady
$ws1.Range("D3").Value = "2016-44-21 04:44:51"

# Hyperlinks keep pointing at the same external targets as before; only the
# displayed text needs to track which file now sits in which row.
$ws1Link2 = "https://github.com/OpenLocalizationTest/oltest/blob/8ed93f2b0908d8f1d52f1d70f92f99241a5c0641/e2e/$fileA.md"
$ws1Link3 = "https://github.com/OpenLocalizationTest/oltest/blob/8ed93f2b0908d8f1d52f1d70f92f99241a5c0641/e2e/$fileB.md"

$ws1.Cells.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $ws1Link2, "", "", "$fileB.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), $ws1Link3, "", "", "$fileA.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "$fileB.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = $statusHandedBack
$ws2.Range("D2").Value = "$fileB.cdf4e48b29d753218640c23b7276de04b8bdaf88.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-21 04:43:49"
$ws2.Range("F2").Value = "$fileB.md"
$ws2.Range("G2").Value = "$fileB.cdf4e48b29d753218640c23b7276de04b8bdaf88.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-03-21 04:44:15"
$ws2.Range("I2").Value = "Include"

$ws2.Range("A3").Value = "$fileA.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = $statusReady
$ws2.Range("D3").Value = "$fileA.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-21 04:44:48"
$ws2.Range("F3").Value = "$fileA.md"
$ws2.Range("G3").Value = "$fileA.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-03-21 04:44:15"
$ws2.Range("I3").Value = "Include"

$ws2LinkA2 = "https://github.com/OpenLocalizationTest/oltest/blob/8ed93f2b0908d8f1d52f1d70f92f99241a5c0641/e2e/$fileA.md"
$ws2LinkB2 = "https://github.com/OpenLocalizationTest/oltest/blob/8ed93f2b0908d8f1d52f1d70f92f99241a5c0641/e2e/$fileA.md"
$ws2LinkD2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2f6946f5fb8b38debaf6ab0f0bd58056723d6c3d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$fileA.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.zh-cn.xlf"
$ws2LinkF2 = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b9df5aff37938bc3165d0b0e1ba9e5e34158b9c5/e2e/$fileA.md"
$ws2LinkG2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4da0a5c1482191a88e8e3b162bd99e2100358a3d/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$fileA.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.zh-cn.xlf"

$ws2LinkA3 = "https://github.com/OpenLocalizationTest/oltest/blob/8ed93f2b0908d8f1d52f1d70f92f99241a5c0641/e2e/$fileB.md"
$ws2LinkB3 = "https://github.com/OpenLocalizationTest/oltest/blob/8ed93f2b0908d8f1d52f1d70f92f99241a5c0641/e2e/$fileB.md"
$ws2LinkD3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2f6946f5fb8b38debaf6ab0f0bd58056723d6c3d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$fileB.cdf4e48b29d753218640c23b7276de04b8bdaf88.zh-cn.xlf"
$ws2LinkF3 = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b9df5aff37938bc3165d0b0e1ba9e5e34158b9c5/e2e/$fileB.md"
$ws2LinkG3 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4da0a5c1482191a88e8e3b162bd99e2100358a3d/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$fileB.cdf4e48b29d753218640c23b7276de04b8bdaf88.zh-cn.xlf"

$ws2.Cells.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $ws2LinkA2, "", "", "$fileB.md")
$ws2.Hyperlinks.Add($ws2.Range("B2"), $ws2LinkB2, "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), $ws2LinkD2, "", "", "$fileB.cdf4e48b29d753218640c23b7276de04b8bdaf88.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), $ws2LinkF2, "", "", "$fileB.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), $ws2LinkG2, "", "", "$fileB.cdf4e48b29d753218640c23b7276de04b8bdaf88.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), $ws2LinkA3, "", "", "$fileA.md")
$ws2.Hyperlinks.Add($ws2.Range("B3"), $ws2LinkB3, "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), $ws2LinkD3, "", "", "$fileA.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F3"), $ws2LinkF3, "", "", "$fileA.md")
$ws2.Hyperlinks.Add($ws2.Range("G3"), $ws2LinkG3, "", "", "$fileA.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.zh-cn.xlf")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "$fileB.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = $statusHandedBack
$ws3.Range("D2").Value = "$fileB.cdf4e48b29d753218640c23b7276de04b8bdaf88.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-21 04:43:53"
$ws3.Range("F2").Value = "$fileB.md"
$ws3.Range("G2").Value = "$fileB.cdf4e48b29d753218640c23b7276de04b8bdaf88.de-de.xlf"
$ws3.Range("H2").Value = "2016-03-21 04:44:22"
$ws3.Range("I2").Value = "Include"

$ws3.Range("A3").Value = "$fileA.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = $statusReady
$ws3.Range("D3").Value = "$fileA.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-21 04:44:51"
$ws3.Range("F3").Value = "$fileA.md"
$ws3.Range("G3").Value = "$fileA.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.de-de.xlf"
$ws3.Range("H3").Value = "2016-03-21 04:44:22"
$ws3.Range("I3").Value = "Include"

$ws3LinkA2 = "https://github.com/OpenLocalizationTest/oltest/blob/8ed93f2b0908d8f1d52f1d70f92f99241a5c0641/e2e/$fileA.md"
$ws3LinkB2 = "https://github.com/OpenLocalizationTest/oltest/blob/8ed93f2b0908d8f1d52f1d70f92f99241a5c0641/e2e/$fileA.md"
$ws3LinkD2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5d6ebebd414ff00125800fd1d53f6a2014942daa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$fileA.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.de-de.xlf"
$ws3LinkF2 = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/aabc32508b60bf5f97fd306775e4ffac24a4dd76/e2e/$fileA.md"
$ws3LinkG2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/408141affcc5a2f99a8f5e24c28b7f73d0864cee/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$fileA.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.de-de.xlf"

$ws3LinkA3 = "https://github.com/OpenLocalizationTest/oltest/blob/8ed93f2b0908d8f1d52f1d70f92f99241a5c0641/e2e/$fileB.md"
$ws3LinkB3 = "https://github.com/OpenLocalizationTest/oltest/blob/8ed93f2b0908d8f1d52f1d70f92f99241a5c0641/e2e/$fileB.md"
$ws3LinkD3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5d6ebebd414ff00125800fd1d53f6a2014942daa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$fileB.cdf4e48b29d753218640c23b7276de04b8bdaf88.de-de.xlf"
$ws3LinkF3 = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/aabc32508b60bf5f97fd306775e4ffac24a4dd76/e2e/$fileB.md"
$ws3LinkG3 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/408141affcc5a2f99a8f5e24c28b7f73d0864cee/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$fileB.cdf4e48b29d753218640c23b7276de04b8bdaf88.de-de.xlf"

$ws3.Cells.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $ws3LinkA2, "", "", "$fileB.md")
$ws3.Hyperlinks.Add($ws3.Range("B2"), $ws3LinkB2, "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), $ws3LinkD2, "", "", "$fileB.cdf4e48b29d753218640c23b7276de04b8bdaf88.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), $ws3LinkF2, "", "", "$fileB.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), $ws3LinkG2, "", "", "$fileB.cdf4e48b29d753218640c23b7276de04b8bdaf88.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), $ws3LinkA3, "", "", "$fileA.md")
$ws3.Hyperlinks.Add($ws3.Range("B3"), $ws3LinkB3, "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), $ws3LinkD3, "", "", "$fileA.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F3"), $ws3LinkF3, "", "", "$fileA.md")
$ws3.Hyperlinks.Add($ws3.Range("G3"), $ws3LinkG3, "", "", "$fileA.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.de-de.xlf")
